# Remove the post "「デザイン」タスミーム" (row 705) from the sheet.
# All subsequent rows shift up by one, which matches the updated dimension
# (A1:C818 -> A1:C817).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(705).Delete()
